$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$shp = $hdr.Range.InlineShapes.Item(1)
$shp.Select()
$sel = $word.Selection
try {
    Write-Host ("ShapeRange count=" + $sel.ShapeRange.Count)
} catch { Write-Host ("err: " + $_) }
